# Update the "approval"/"threshold"/etc. boolean matrix (columns B:G, rows 2-15)
# so that the satisfaction logic that used to live in separate functions is now
# combined into a single function - which changed which boolean flags get set
# for each project row.
#
# Final (after) state of columns B,C,D,E,F,G for rows 2-15:
#   row  B  C  D  E  F  G
#    2   1  1  0  0  0  1
#    3   1  1  0  1  0  0
#    4   1  1  0  0  0  1
#    5   1  1  0  0  0  1
#    6   1  1  0  0  0  1
#    7   1  1  0  0  0  1
#    8   1  1  0  1  0  0
#    9   1  1  0  0  0  1
#   10   1  1  0  1  0  0
#   11   1  1  0  1  0  0
#   12   1  1  0  0  0  1
#   13   1  1  0  1  0  0
#   14   1  1  0  1  0  0
#   15   1  1  0  1  0  0

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @{
    2  = @($true,  $true, $false, $false, $false, $true)
    3  = @($true,  $true, $false, $true,  $false, $false)
    4  = @($true,  $true, $false, $false, $false, $true)
    5  = @($true,  $true, $false, $false, $false, $true)
    6  = @($true,  $true, $false, $false, $false, $true)
    7  = @($true,  $true, $false, $false, $false, $true)
    8  = @($true,  $true, $false, $true,  $false, $false)
    9  = @($true,  $true, $false, $false, $false, $true)
    10 = @($true,  $true, $false, $true,  $false, $false)
    11 = @($true,  $true, $false, $true,  $false, $false)
    12 = @($true,  $true, $false, $false, $false, $true)
    13 = @($true,  $true, $false, $true,  $false, $false)
    14 = @($true,  $true, $false, $true,  $false, $false)
    15 = @($true,  $true, $false, $true,  $false, $false)
}

$cols = @("B", "C", "D", "E", "F", "G")

foreach ($r in $rows.Keys) {
    $values = $rows[$r]
    for ($i = 0; $i -lt $cols.Length; $i++) {
        $ws.Range("$($cols[$i])$r").Value = $values[$i]
    }
}
